# Fixing secondary phases to work with CRR Code
#
# The header row (row 1) gains a new "rays_present" column right after
# "filename", and three new ratio columns are inserted right after
# "Diad2_prom/std_betweendiads". Concretely this shifts every header from
# column C onward one slot to the right (C1 <- old C1's neighbour set),
# with the very last four original headers (Diad1_Median_Bck,
# Diad2_Median_Bck, C13_HB2_abs_prom_ratio, Diad2_HB2_Valley_prom) ending
# up appended as four brand-new trailing columns AM1:AP1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 ("filename") is unchanged.
$ws.Range("C1").Value  = 'rays_present'
$ws.Range("D1").Value  = 'approx_split'
$ws.Range("E1").Value  = 'Diad1_pos'
$ws.Range("F1").Value  = 'Diad2_pos'
$ws.Range("G1").Value  = 'HB1_pos'
$ws.Range("H1").Value  = 'HB2_pos'
$ws.Range("I1").Value  = 'C13_pos'
$ws.Range("J1").Value  = 'Diad1_abs_prom'
$ws.Range("K1").Value  = 'Diad2_abs_prom'
$ws.Range("L1").Value  = 'HB1_abs_prom'
$ws.Range("M1").Value  = 'HB2_abs_prom'
$ws.Range("N1").Value  = 'C13_abs_prom'
$ws.Range("O1").Value  = 'Mean_abs_HB_prom'
$ws.Range("P1").Value  = 'Diad2_HB2_abs_prom_ratio'
$ws.Range("Q1").Value  = 'Diad1_HB1_abs_prom_ratio'
$ws.Range("R1").Value  = 'Diad1_rel_prom'
$ws.Range("S1").Value  = 'Diad2_rel_prom'
$ws.Range("T1").Value  = 'HB1_rel_prom'
$ws.Range("U1").Value  = 'HB2_rel_prom'
$ws.Range("V1").Value  = 'C13_rel_prom'
$ws.Range("W1").Value  = 'Diad1_HB1_abs_prom_ratio'
$ws.Range("X1").Value  = 'Diad2_HB2_abs_prom_ratio'
$ws.Range("Y1").Value  = 'Diad1_HB1_Valley_prom'
$ws.Range("Z1").Value  = 'Diad2_HB2_abs_prom_ratio'
$ws.Range("AA1").Value = 'Mean_Diad_HB_Valley_prom'
$ws.Range("AB1").Value = 'Mean_abs_HB_prom'
$ws.Range("AC1").Value = 'Diad1_prom/std_betweendiads'
$ws.Range("AD1").Value = 'Diad2_prom/std_betweendiads'
$ws.Range("AE1").Value = 'Av_Diad_prom/std_betweendiads'
$ws.Range("AF1").Value = 'C13_prom/HB2_prom'
$ws.Range("AG1").Value = 'Av_Diad_HB_prom_ratio'
$ws.Range("AH1").Value = 'Diad2_height'
$ws.Range("AI1").Value = 'HB2_height'
$ws.Range("AJ1").Value = 'C13_height'
$ws.Range("AK1").Value = 'Diad1_height'
$ws.Range("AL1").Value = 'HB1_height'

# Four brand-new trailing header cells. Copy AL1's format (bold header
# font + border + centered/top alignment) across before overwriting the
# text, so the new cells keep the same "s=1" header style as the rest of
# row 1 instead of defaulting to no style.
$ws.Range("AL1").Copy()

$ws.Range("AM1").PasteSpecial(-4122)
$ws.Range("AM1").Value = 'Diad1_Median_Bck'

$ws.Range("AN1").PasteSpecial(-4122)
$ws.Range("AN1").Value = 'Diad2_Median_Bck'

$ws.Range("AO1").PasteSpecial(-4122)
$ws.Range("AO1").Value = 'C13_HB2_abs_prom_ratio'

$ws.Range("AP1").PasteSpecial(-4122)
$ws.Range("AP1").Value = 'Diad2_HB2_Valley_prom'
